$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (id) and C (speaker_variant) for rows 2-14.
# Column D (is_prefered) is cleared for every row (no more "x" markers).
# Note: C2's target text starts with a literal apostrophe ('tgeru). When
# assigned through .Value, Excel treats a single leading apostrophe as the
# "force text" prefix and strips it, so it must be escaped by doubling it
# (''tgeru) to end up with a single stored apostrophe.
$rows = @(
    @{ Row = 2;  B = "#'tgeru";     C = "''tgeru" }
    @{ Row = 3;  B = "#achis";      C = "achis" }
    @{ Row = 4;  B = "#golia";      C = "golia" }
    @{ Row = 5;  B = "#ionath";     C = "ionath" }
    @{ Row = 6;  B = "#saul";       C = "saul" }
    @{ Row = 7;  B = "#1.-philist"; C = "1. philist" }
    @{ Row = 8;  B = "#david";      C = "david" }
    @{ Row = 9;  B = "#2.-philist"; C = "2. philist" }
    @{ Row = 10; B = "#eliab";      C = "eliab" }
    @{ Row = 11; B = "#schild";     C = "schild" }
    @{ Row = 12; B = "#isai";       C = "isai" }
    @{ Row = 13; B = "#abner";      C = "abner" }
    @{ Row = 14; B = "#goliat";     C = "goliat" }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = ""
}
